$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.622.01"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "2.572.82"
$ws.Range("E3").Value = "  -2.03%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'516.54"
$ws.Range("E5").Value = "  -2.42%  "
$ws.Range("D6").Value = "'138.28"
$ws.Range("E6").Value = "  -4.30%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'0.563"
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("D9").Value = "2.587.47"
$ws.Range("E9").Value = "  -1.87%  "
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").Value = "'0.1000"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("D12").Value = "'0.330"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("D14").Value = "3.028.79"
$ws.Range("E14").Value = "  -2.14%  "
$ws.Range("D15").Value = "58.537.59"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Value = "'20.23"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.562.84"
$ws.Range("E17").Value = "  -2.98%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000131"
$ws.Range("E18").Value = "  -2.56%  "
$ws.Range("D19").Value = "'336.08"
$ws.Range("E19").Value = "  -1.53%  "
$ws.Range("E20").Value = "  -2.43%  "
$ws.Range("E21").Value = "  -3.09%  "
$ws.Range("D22").Value = "'6.41"
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'65.88"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").Value = "'0.166"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("D27").Value = "'0.996"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").Value = "'6.99"
$ws.Range("E28").Value = "  -1.26%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -8.96%  "
$ws.Range("D31").Value = "'5.91"
$ws.Range("E31").Value = "  -5.32%  "
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("E33").Value = "  -3.05%  "
$ws.Range("D34").Value = "'148.43"
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").Value = "'3.92"
$ws.Range("E35").Value = "  -4.11%  "
$ws.Range("E36").Value = "  -2.93%  "
$ws.Range("D37").Value = "'36.14"
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("D38").Value = "'0.833"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").Value = "'0.818"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "'271.20"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("D44").Value = "'10.73"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.584"
$ws.Range("E45").Value = "  -1.58%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.0944"
$ws.Range("E46").Value = "  -2.67%  "
$ws.Range("D47").Value = "'0.0516"
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("D48").Value = "1.971.98"
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("D49").Value = "'18.33"
$ws.Range("E49").Value = "  -3.18%  "
$ws.Range("E50").Value = "  -3.94%  "
$ws.Range("E51").Value = "  -2.96%  "
